$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "MM2316800126"
$ws.Range("A4").Value = "MM2316800148"
